$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.012.21'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.53%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.695.58'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.11%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.02%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '649.83'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -4.05%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '162.38'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.92%  '

# Row 7
$ws.Range('E7').Value = '  -0.02%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.502'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +1.18%  '

# Row 9
$ws.Range('E9').Value = '  -1.01%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.23'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.80%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.444'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.54%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000233'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.54%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.311.83'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.19%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.84'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.60%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.692.61'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.04%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '69.930.38'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.41%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.117'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.62%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '16.02'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.18%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.53'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.84%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.41'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +5.93%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '471.88'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.23%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.655'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.41%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '80.21'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.42%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.840.73'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.12%  '

# Row 25
$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000128'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.84%  '

# Row 26
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.03%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.01'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.82%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.16'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.40%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.66'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.53%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.72'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.79%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.169'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +4.13%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.02'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.19%  '

# Row 33
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.55'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.81%  '

# Row 34
$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.34%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '26.80'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.50%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.689.54'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.05%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.47'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.14%  '

# Row 38
$ws.Range('E38').Value = '  -0.10%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '180.42'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +7.89%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.91'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -5.24%  '

# Row 41
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.23'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.51%  '

# Row 42
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.998'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.15%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0905'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.31%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.934'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.09%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.88'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +3.74%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '29.53'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +4.53%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '46.61'
$ws.Range('D47').ClearFormats()

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.000275'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.22%  '

# Row 49
$ws.Range('E49').Value = '  -3.37%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.88'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.08%  '

# Row 51
$ws.Range('E51').Value = '  -3.89%  '
